$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r0")

$rows = 30,31,32,33
$names = "ExoT_r0_script_12v0","ExoT_r0_script_12v1","ExoT_r0_script_12v2","ExoT_r0_script_12v3"
$flowRates = "2 mL/hr","5 mL/hr","10 mL/hr","25 mL/hr"

for ($i = 0; $i -lt 4; $i++) {
  $ws.Cells.Item($rows[$i], 1).Value = $names[$i]
}
for ($i = 0; $i -lt 4; $i++) {
  $r = $rows[$i]
  $ws.Cells.Item($r, 2).Value = "Andrew's flow rate optimization"
  $ws.Cells.Item($r, 3).Value = "5 mL"
  $ws.Cells.Item($r, 4).Value = "5 mL"
  $ws.Cells.Item($r, 5).Value = "1 hour"
  $ws.Cells.Item($r, 6).Value = "1 mL"
  $ws.Cells.Item($r, 7).Value = $flowRates[$i]
  $ws.Cells.Item($r, 8).Value = "15 mL/hr"
  $ws.Cells.Item($r, 9).Value = "200-800-1000"
  $ws.Cells.Item($r, 10).Value = "2 mins"
  $ws.Cells.Item($r, 11).Value = "N"
}

$ws.Activate()
$ws.Range("A30").Select()
